$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5874.75
$ws.Range("I76").Value = 5833
$ws.Range("K76").Value = 5833
$ws.Range("M76").Value = -5518
$ws.Range("H79").Value = 5874.75
$ws.Range("I79").Value = 5833
$ws.Range("K79").Value = 5833
$ws.Range("M79").Value = -4741
$ws.Range("H104").Value = 388.4
$ws.Range("I104").Value = 418.22223
$ws.Range("K104").Value = 1254.66669
$ws.Range("M104").Value = 492.33331
$ws.Range("H116").Value = 4892.6
$ws.Range("I116").Value = 4398
$ws.Range("K116").Value = 4398
$ws.Range("M116").Value = -956
$ws.Range("H131").Value = 3226
$ws.Range("I131").Value = 2615.7646
$ws.Range("J131").Value = 5300.8
$ws.Range("K131").Value = 7847.293799999999
$ws.Range("L131").Value = 15902.4
$ws.Range("M131").Value = -2807.293799999999
$ws.Range("N131").Value = -25982.4
$ws.Range("H132").Value = 47623616
$ws.Range("I132").Value = 71434160
$ws.Range("J132").Value = 2538.4285
$ws.Range("K132").Value = 214302480
$ws.Range("L132").Value = 7615.2855
$ws.Range("M132").Value = -214299950
$ws.Range("N132").Value = -12675.2855
$ws.Range("H138").Value = 2821.4119
$ws.Range("I138").Value = 1234.8572
$ws.Range("J138").Value = 4752.8696
$ws.Range("K138").Value = 3704.5716
$ws.Range("L138").Value = 14258.6088
$ws.Range("M138").Value = 1435.4284
$ws.Range("N138").Value = -24538.6088
$ws.Range("H141").Value = 1477.9565
$ws.Range("I141").Value = 1452.0476
$ws.Range("K141").Value = 4356.142800000001
$ws.Range("M141").Value = 823.8571999999995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3378.0605
$ws.Range("I32").Value = 2033.3455
$ws.Range("K32").Value = 2033.3455
$ws.Range("M32").Value = -1746.3455
$ws.Range("H102").Value = 2978780.8
$ws.Range("I102").Value = 3969806.5
$ws.Range("K102").Value = 3969806.5
$ws.Range("M102").Value = -3968184.5
$ws.Range("H132").Value = 3374.5
$ws.Range("I132").Value = 3501
$ws.Range("J132").Value = 3248
$ws.Range("K132").Value = 10503
$ws.Range("L132").Value = 9744
$ws.Range("M132").Value = -7973
$ws.Range("N132").Value = -14804
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12808
$ws.Range("I31").Value = 1765.5807
$ws.Range("J31").Value = 17917.18
$ws.Range("K31").Value = 1765.5807
$ws.Range("L31").Value = 17917.18
$ws.Range("M31").Value = -1470.5807
$ws.Range("N31").Value = -18507.18
$ws.Range("H34").Value = 12808
$ws.Range("I34").Value = 1765.5807
$ws.Range("J34").Value = 17917.18
$ws.Range("K34").Value = 1765.5807
$ws.Range("L34").Value = 17917.18
$ws.Range("M34").Value = -1563.5807
$ws.Range("N34").Value = -18321.18
$ws.Range("H58").Value = 2750.9565
$ws.Range("J58").Value = 3451.5557
$ws.Range("L58").Value = 3451.5557
$ws.Range("N58").Value = -3857.5557
$ws.Range("H132").Value = 106854.445
$ws.Range("I132").Value = 64593.938
$ws.Range("J132").Value = 444938.5
$ws.Range("K132").Value = 193781.814
$ws.Range("L132").Value = 1334815.5
$ws.Range("M132").Value = -191251.814
$ws.Range("N132").Value = -1339875.5
$ws.Range("H136").Value = 2750.9565
$ws.Range("J136").Value = 3451.5557
$ws.Range("L136").Value = 10354.6671
$ws.Range("N136").Value = -15454.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 7270.7144
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 8315.833000000001
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 24947.499
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -25203.499
$ws.Range("H92").Value = 866
$ws.Range("I92").Value = 409.66666
$ws.Range("J92").Value = 1037.125
$ws.Range("K92").Value = 1228.99998
$ws.Range("L92").Value = 3111.375
$ws.Range("M92").Value = 19.00001999999995
$ws.Range("N92").Value = -5607.375
$ws.Range("H93").Value = 5000
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H96").Value = 14090
$ws.Range("H97").Value = 391
$ws.Range("I97").Value = 367
$ws.Range("J97").Value = 410.2
$ws.Range("K97").Value = 1101
$ws.Range("L97").Value = 1230.6
$ws.Range("M97").Value = -605
$ws.Range("N97").Value = -2222.6
$ws.Range("H98").Value = 509.375
$ws.Range("J98").Value = 682
$ws.Range("L98").Value = 2046
$ws.Range("N98").Value = -5042
$ws.Range("H99").Value = 1571.5
$ws.Range("I99").Value = 1571.5
$ws.Range("K99").Value = 4714.5
$ws.Range("M99").Value = -2468.5
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2189
$ws.Range("H102").Value = 13977.272
$ws.Range("H103").Value = 100
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 140699
$ws.Range("J133").Value = 140699
$ws.Range("L133").Value = 140699
$ws.Range("N133").Value = -150819
$ws.Range("H136").Value = 42000
$ws.Range("I136").Value = 30000
$ws.Range("J136").Value = 48000
$ws.Range("K136").Value = 90000
$ws.Range("L136").Value = 144000
$ws.Range("M136").Value = -87450
$ws.Range("N136").Value = -149100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11114503
$ws.Range("I61").Value = 11114503
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 11114503
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -11114301
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 11114503
$ws.Range("I113").Value = 11114503
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 11114503
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -11112333
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 11358.462
$ws.Range("I132").Value = 12712.223
$ws.Range("K132").Value = 38136.669
$ws.Range("M132").Value = -35606.669
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1883.2778
$ws.Range("I122").Value = 1460.2667
$ws.Range("K122").Value = 4380.800099999999
$ws.Range("M122").Value = -1930.800099999999
$ws.Range("H132").Value = 37446050
$ws.Range("I132").Value = 41672560
$ws.Range("K132").Value = 125017680
$ws.Range("M132").Value = -125015150
